$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "legenda"
$ws.Range("B1").Value = "area"
$ws.Range("D1").Value = "area_km2"
$ws.Range("B2").Value = 158792.012421
